# household_member.xlsx update:
#  - reword 3 prompts to reference {{member_name}}
#  - insert a new "note" row (age odd/even example) into the survey sheet
#  - add a "calculates" sheet with a calculation_name/calculation example

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. survey sheet: reword existing prompts to reference {{member_name}}
# ---------------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

$survey.Range("D4").Value = "Enter age of {{member_name}}:"
$survey.Range("D5").Value = "Enter sex of {{member_name}}:"

# insert the new blank "note" row above income_contribution first, so the
# row/values below can be filled in afterwards in the same order the
# original author used
$survey.Rows.Item(6).Insert()
$survey.Rows.Item(6).RowHeight = 31

# the income_contribution row is now row 7 - reword its prompt
$survey.Cells.Item(7, 4).Value = "Does {{member_name}} contribute to the household income?"

# ---------------------------------------------------------------------------
# 2. choices sheet: move the selection (cosmetic, matches authored file)
# ---------------------------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")
$choices.Columns.Item(5).Select()

# ---------------------------------------------------------------------------
# 3. add the "calculates" sheet after "model"
# ---------------------------------------------------------------------------
$model = $wb.Worksheets.Item("model")
$calculates = $wb.Worksheets.Add($null, $model)
$calculates.Name = "calculates"

$calculates.Columns.Item(1).ColumnWidth = 19.453125
$calculates.Columns.Item(2).ColumnWidth = 44.7265625

$calculates.Cells.Item(1, 1).Value = "calculation_name"
$calculates.Cells.Item(1, 2).Value = "calculation"
$calculates.Rows.Item(1).RowHeight = 13.5

# ---------------------------------------------------------------------------
# 4. survey sheet: fill in the "note" row added above
# ---------------------------------------------------------------------------
$survey.Cells.Item(6, 1).Value = "note"

# ---------------------------------------------------------------------------
# 5. calculates sheet: the actual calculation referenced by the note above
# ---------------------------------------------------------------------------
$calculates.Cells.Item(2, 1).Value = "ageIsOddOrEven"
$calculates.Cells.Item(2, 2).Value = "((data('age') % 2) == 1) ? ""odd"" : ""even"""

$calculates.Range("A2").Select()

# ---------------------------------------------------------------------------
# 6. survey sheet: the note text evaluating the calculation above
# ---------------------------------------------------------------------------
$survey.Cells.Item(6, 4).Value = "{{member_name}} age is {{evaluate calculates.ageIsOddOrEven}} in {{setting 'table_id'}} for {{metadata 'instanceName'}}"

# ---------------------------------------------------------------------------
# 7. leave the focus back on the survey sheet, like the authored workbook
# ---------------------------------------------------------------------------
$survey.Range("D7").Select()
